$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace values in B1:E3 with the values currently found in O1:R3 and AN1:AQ3
# (columns B:C take the O/R values, columns D:E take the AN/AQ values)
$ws.Range("B1").Value = $ws.Range("O1").Value()
$ws.Range("C1").Value = $ws.Range("R1").Value()
$ws.Range("D1").Value = $ws.Range("AN1").Value()
$ws.Range("E1").Value = $ws.Range("AQ1").Value()

$ws.Range("B2").Value = $ws.Range("O2").Value()
$ws.Range("C2").Value = $ws.Range("R2").Value()
$ws.Range("D2").Value = $ws.Range("AN2").Value()
$ws.Range("E2").Value = $ws.Range("AQ2").Value()

$ws.Range("B3").Value = $ws.Range("O3").Value()
$ws.Range("C3").Value = $ws.Range("R3").Value()
$ws.Range("D3").Value = $ws.Range("AN3").Value()
$ws.Range("E3").Value = $ws.Range("AQ3").Value()

# Update the selected range to match the new selection B1:E3
$null = $ws.Range("B1:E3").Select()
